$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7 (between "Chiapas" and "Coahuila"), shifting
# subsequent rows down, then set its value to the new state name.
$ws.Rows("7:7").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)
$ws.Range("A7").Value = "Ciudad de México"

# Match the active cell/selection shown in the saved workbook.
$ws.Range("A8").Select()
